$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Regenerated the handoff report: status text flips from the old
# "handed back" state to "ready for handoff", and the generation
# timestamps move forward to the new run.
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

$wsOverview.Range("G2").Value = "2016-08-23 11:01:22"
$wsZhCn.Range("H2").Value = "2016-08-23 11:01:17"
$wsDeDe.Range("H2").Value = "2016-08-23 11:01:22"

# The shorter status text needs less room, so narrow the columns that
# display it (~16.33 chars -> engine rounds to the pixel grid closest
# to the target 17.22-wide column used by the authoring tool).
$wsOverview.Columns("E").ColumnWidth = 16.33
$wsOverview.Columns("F").ColumnWidth = 16.33
$wsZhCn.Columns("C").ColumnWidth = 16.33
$wsDeDe.Columns("C").ColumnWidth = 16.33
